$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.130.36'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.873.98'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.59%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.84'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9992'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5057'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3845'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09046'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.121'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.67'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.92%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.368'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.30%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.885.30'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.278'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.35'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06649'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.19'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9990'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.136'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.156.53'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.47'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.262'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.548'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.091.90'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.72%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.84'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '157.07'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.92'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1064'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.19%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.615'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.594'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.467'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06598'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02410'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.71%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2195'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.293'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.213'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6409'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.51'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.928'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9988'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.30'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.27%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6037'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.276'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.236'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.006'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '121.39'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.40%  '
